$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing data row (47) down through the
# new rows (48-90) so the new records match the sheet style (s="2").
$ws.Range("A47:G47").Copy()
$ws.Range("A48:G90").PasteSpecial(-4122)

# New player records (Barcelona squad, then Paris Saint-Germain squad).
# Columns: Name|Position|Age|Country|Squad number|Club|Height
$newPlayers = @(
  "Marc-Andre ter Stegen|Goalkeeper|22|Germany|1|Barcelona|189",
  "Claudio Bravo|Goalkeeper|31|Spain|13|Barcelona|185",
  "Martin Montoya|Defender|23|Spain|2|Barcelona|174",
  "Gerard Pique|Defender|27|Spain|3|Barcelona|192",
  "Javier Mascherano|Defender|30|Argentina|14|Barcelona|174",
  "Marc Bartra|Defender|23|Spain|15|Barcelona|183",
  "Douglas|Defender|24|Brazil|16|Barcelona|170",
  "Jordi Alba|Defender|25|Spain|18|Barcelona|170",
  "Adriano|Defender|30|Brazil|21|Barcelona|172",
  "Daniel Alves|Defender|31|Brazil|22|Barcelona|173",
  "Thomas Vermaelen|Defender|29|Belgium|23|Barcelona|180",
  "Jeremy Mathieu|Defender|31|France|24|Barcelona|192",
  "Ivan Rakitic|Midfield|26|Croatia|4|Barcelona|184",
  "Sergio Busquets|Midfield|26|Spain|5|Barcelona|189",
  "Xavi Hernandez|Midfield|34|Spain|6|Barcelona|170",
  "Andres Iniesta|Midfield|30|Spain|8|Barcelona|170",
  "Rafinha|Midfield|21|Spain|12|Barcelona|174",
  "Sergi Roberto|Midfield|22|Spain|20|Barcelona|177",
  "Pedro Rodríguez|Forward|27|Spain|7|Barcelona|169",
  "Luis Suarez|Forward|27|Uruguay|9|Barcelona|181",
  "Lionel Messi|Forward|27|Argentina|10|Barcelona|169",
  "Neymar|Forward|22|Brazil|11|Barcelona|175",
  "Munir El Haddadi|Forward|19|Spain|31|Barcelona|175",
  "Sandro Ramirez|Forward|19|Spain|29|Barcelona|175",
  "Salvatore Sirigu|Goalkeeper|27|Italy|30|Paris|192",
  "Thiago Silva|Defender|30|Brazil|2|Paris|183",
  "Marquinhos|Defender|20|Brazil|5|Paris|183",
  "Maxwell|Defender|33|Brazil|17|Paris|177",
  "Lucas Digne|Defender|21|France|21|Paris|178",
  "Gregory van der Wiel|Defender|26|Netherlands|23|Paris|172",
  "David Luiz|Defender|27|Brazil|32|Paris|188",
  "Yohan Cabaye|Midfield|28|France|4|Paris|175",
  "Lucas|Midfield|22|Brazil|7|Paris|172",
  "Thiago Motta|Midfield|32|Italy|8|Paris|183",
  "Blaise Matuidi|Midfield|27|France|14|Paris|175",
  "Clement Chantome|Midfield|27|France|20|Paris|180",
  "Marco Verratti|Midfield|22|Italy|24|Paris|165",
  "Adrien Rabiot|Midfield|19|France|25|Paris|188",
  "Javier Pastore|Midfield|25|Argentina|27|Paris|187",
  "Edinson Cavani|Forward|27|Uruguay|9|Paris|184",
  "Zlatan Ibrahimovic|Forward|33|Sweden|10|Paris|194",
  "Jean-Christophe Bahebeck|Forward|21|France|15|Paris|182",
  "Ezequiel Lavezzi|Forward|29|Argentina|22|Paris|173"
)

$startRow = 48
for ($i = 0; $i -lt $newPlayers.Count; $i++) {
  $fields = $newPlayers[$i].Split("|")
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $fields[0]
  $ws.Cells.Item($r, 2).Value = $fields[1]
  $ws.Cells.Item($r, 3).Value = [int]$fields[2]
  $ws.Cells.Item($r, 4).Value = $fields[3]
  $ws.Cells.Item($r, 5).Value = [int]$fields[4]
  $ws.Cells.Item($r, 6).Value = $fields[5]
  $ws.Cells.Item($r, 7).Value = [int]$fields[6]
}

# Restore the active cell/selection to match the edited workbook state.
$ws.Range("C85").Select()
